$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "07/27/2025"
$ws.Range("B35").Value = 0.0004206000000000001
$ws.Range("C35").Value = 118877.7936281502
$ws.Range("D35").Value = 50
